$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cells (data period label and download timestamp)
$ws.Cells.Item(4, 2).Value = "2024-09-01 ~ 2024-09-30"
$ws.Cells.Item(5, 2).Value = "2024년 10월 02일 19시 45분 47초"

# The data table shrank from 31 days (Aug) to 30 days (Sep); drop the last
# data row (old row 38) so the trailing blank row shifts up to row 38.
$ws.Rows(38).Delete()

# Force columns A (dates) and C:F (counts) to be stored as literal text,
# matching the original inline-string cell type (otherwise the COM layer
# would coerce "2024-09-30"-like strings to dates and "28"-like strings
# to numbers).
$ws.Range("A8:A37").NumberFormat = "@"
$ws.Range("C8:F37").NumberFormat = "@"

# Replace the 30-day August dataset with the new September dataset.
$data = @(
    @("2024-09-30", "월", "28", "1", "0", "27"),
    @("2024-09-29", "일", "13", "0", "0", "13"),
    @("2024-09-28", "토", "9", "0", "0", "9"),
    @("2024-09-27", "금", "21", "0", "0", "21"),
    @("2024-09-26", "목", "16", "0", "0", "16"),
    @("2024-09-25", "수", "21", "0", "0", "21"),
    @("2024-09-24", "화", "16", "0", "0", "16"),
    @("2024-09-23", "월", "19", "0", "0", "19"),
    @("2024-09-22", "일", "10", "0", "0", "10"),
    @("2024-09-21", "토", "17", "0", "0", "17"),
    @("2024-09-20", "금", "17", "0", "0", "17"),
    @("2024-09-19", "목", "24", "1", "1", "23"),
    @("2024-09-18", "수", "31", "1", "0", "30"),
    @("2024-09-17", "화", "17", "1", "0", "16"),
    @("2024-09-16", "월", "28", "5", "1", "23"),
    @("2024-09-15", "일", "43", "4", "1", "38"),
    @("2024-09-14", "토", "15", "1", "0", "14"),
    @("2024-09-13", "금", "26", "0", "0", "26"),
    @("2024-09-12", "목", "12", "0", "0", "12"),
    @("2024-09-11", "수", "17", "0", "0", "17"),
    @("2024-09-10", "화", "19", "0", "0", "19"),
    @("2024-09-09", "월", "26", "2", "0", "24"),
    @("2024-09-08", "일", "25", "0", "1", "24"),
    @("2024-09-07", "토", "15", "0", "0", "15"),
    @("2024-09-06", "금", "16", "0", "0", "16"),
    @("2024-09-05", "목", "15", "0", "0", "15"),
    @("2024-09-04", "수", "13", "0", "0", "13"),
    @("2024-09-03", "화", "30", "0", "0", "30"),
    @("2024-09-02", "월", "14", "0", "0", "14"),
    @("2024-09-01", "일", "14", "0", "0", "14")
)

$row = 8
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
    $row++
}
